$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARN: replace failed for:" $old
    }
}

# --- Title ---
Replace-Text "Beyond Microcosm: Unraveling Life's Enigmatic Foundation" "Art History: Exploring the Evolution of Visual Expressions"

# --- Author name (merges "Dr" + "." + " Artemis Augustine" runs into one) ---
Replace-Text "Dr. Artemis Augustine" "Moira Richards"

# --- Email (two separate edits, keep the middle "." run intact) ---
Replace-Text "artemis" "rmoira@eduworld"
Replace-Text "augustine@astralink.com" "school"

# --- Body paragraph 1 ---
Replace-Text "Within the intricate tapestry of life, the mysteries of subatomic particles dance in an enigmatic symphony of creation" "Art history is a journey through time, a visual narrative that reflects the evolution of human society"
Replace-Text " Beyond the realm of visible matter, a microcosm of invisible forces orchestrates the very essence of existence" " It invites us to explore the depths of human imagination and creativity, showcasing the diverse expressions of countless artists who have gifted us with masterpieces that transcend time and boundaries"
Replace-Text " As we delve into the uncharted depths of the subatomic realm, we encounter quarks, bosons, and leptons - the elementary particles that form the building blocks of our universe" " From the intricate cave paintings of our ancestors to the vibrant street art of today, art history is a testament to the resilience of the human spirit"
Replace-Text " The interactions between these particles, governed by the laws of quantum mechanics, shape the properties of matter and energy, giving rise to the captivating phenomena that define our physical world" " Herein, we will delve into this extraordinary tapestry of visual expressions, traversing the ages and cultures to uncover the stories behind the artistry"

Replace-Text "Unraveling the complexities of the subatomic realm has led to groundbreaking discoveries in diverse fields, from physics to biology" "In the realm of art history, we encounter epochs of remarkable artistic achievements, each period characterized by its unique aesthetics, techniques, and cultural influences"
Replace-Text " The exploration of subatomic particles has unlocked the secrets of nuclear energy, enabling the harnessing of immense power and pioneering new avenues for energy production" " The Renaissance, like a blooming flower, witnessed the rebirth of classical ideals in art, exemplified by the works of Da Vinci, Michelangelo, and Raphael"
Replace-Text " In the realm of medicine, the advent of particle accelerators has revolutionized cancer treatment, delivering targeted radiation therapy with remarkable precision" " The Baroque era awed with its extravagance, indulging in elaborate forms and dramatic lighting, while the Impressionists broke free from tradition, capturing the fleeting moments of light and color"
Replace-Text " The study of subatomic particles has also shed light on the fundamental mysteries of life, revealing the intricate mechanisms that govern cellular processes and the genetic code that underlies the diversity of species" " From cave paintings to digital art, the evolution of artistry reveals a kaleidoscope of human ingenuity"

Replace-Text "As we continue to voyage into the uncharted territories of the subatomic realm, we encounter puzzles that challenge our understanding of reality" "Art history not only encompasses Western traditions but also delves into the rich artistic heritage of non-Western cultures"
Replace-Text " The elusive nature of dark matter and dark energy remains a tantalizing enigma, beckoning us to decipher their role in the cosmos" " From the delicate brushstrokes of Chinese calligraphy to the vibrant patterns of African masks, each culture tells its story through its art"
Replace-Text " The quest to unify the four fundamental forces of nature - electromagnetism, weak nuclear force, strong nuclear force, and gravity - drives physicists to seek a comprehensive theory that encompasses the entirety of physical phenomena" " These diverse expressions offer glimpses into worldviews, beliefs, and histories that would otherwise remain hidden"
Replace-Text " With each new discovery, we inch closer to unraveling the profound mysteries that lie at the foundation of life and the universe" " By embracing the global panorama of art, we gain a profound understanding of the human experience across time and space"

# --- Summary paragraph ---
Replace-Text "Our journey into the subatomic realm has unveiled a tapestry of captivating discoveries, transforming our comprehension of the universe and its intricacies" "Art history stands as a testament to the enduring power of human creativity and expression"
Replace-Text " From harnessing nuclear energy and revolutionizing cancer treatment to deciphering the secrets of cellular processes and genetic inheritance, the exploration of subatomic particles has indelibly shaped our understanding of life's fundamental building blocks" " It is a narrative of innovation and perseverance, showcasing the diverse expressions of countless artists throughout history"
Replace-Text " Yet, as we delve deeper into this enigmatic realm, we encounter riddles that tease our intellect and challenge our perception of reality" " By exploring the evolution of art across periods and cultures, we gain insights into the human spirit, its triumphs and struggles, hopes and dreams"
Replace-Text " The nature of dark matter and dark energy remains shrouded in mystery, compelling us to unravel their influence on the cosmos" " Art history not only educates us about aesthetics and techniques but also connects us to the past and present, fostering a deeper understanding of our shared humanity"

# Remove the now-superfluous trailing sentences/runs of the Summary paragraph:
# ". The search for a unified theory ... . As we navigate ..."
Replace-Text ". The search for a unified theory of physics, encompassing all fundamental forces, ignites our curiosity and urges us to seek deeper understanding. As we navigate the uncharted depths of the subatomic realm, we embrace the allure of the unknown, knowing that each revelation brings us closer to comprehending the profound mysteries that underpin the fabric of existence." "."

# --- Trailing empty paragraph ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
